# Update gh-pages to output generated at 456a3b4
# Updates "想去人数" (F) and "最低票价" (G) values on the 展览 and 全部类型
# sheets to reflect freshly scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (rows keyed by their position in that sheet) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 446
$ws1.Range("F3").Value = 25
$ws1.Range("F4").Value = 18
$ws1.Range("F5").Value = 4036
$ws1.Range("F6").Value = 165
$ws1.Range("F7").Value = 51
$ws1.Range("F8").Value = 254
$ws1.Range("G8").Value = 77
$ws1.Range("F9").Value = 25

# --- Sheet "全部类型" (same events, different row offsets) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 446
$ws4.Range("F7").Value = 25
$ws4.Range("F8").Value = 18
$ws4.Range("F9").Value = 4036
$ws4.Range("F10").Value = 165
$ws4.Range("F11").Value = 51
$ws4.Range("F13").Value = 254
$ws4.Range("G13").Value = 77
$ws4.Range("F14").Value = 25
